# Update the "dSF" column (F) values for a set of rows, per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    8  = 4
    10 = -4
    16 = -4
    17 = -3
    20 = -4
    21 = 1
    23 = -17
    27 = 1
    28 = 0
    31 = 4
    36 = 3
    38 = 4
    48 = 0
    52 = -3
    53 = -1
    54 = 1
    62 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
